# Updated cryptos list: apply Price (D) and Volume(1h) (E) cell edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.291.38"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = "'3.100.30"
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'522.97"
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = "'140.84"
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'3.098.67"
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").Value = "'0.438"
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = "'0.384"
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = "'3.631.39"
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").Value = "'26.02"
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = "'57.367.43"
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = "'3.099.28"
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Value = "'337.63"
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").Value = "'66.59"
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("E26").Value = '  -1.45%  '
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").Value = "'0.0₃0911"
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("E32").Value = '  +1.84%  '
$ws.Range("D33").Value = "'20.95"
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("D35").Value = "'156.68"
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = "'27.15"
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").Value = "'0.0658"
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").Value = "'3.94"
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").Value = "'3.140.98"
$ws.Range("D43").Value = "'0.685"
$ws.Range("E43").Value = '  +4.01%  '
$ws.Range("E44").Value = '  +9.84%  '
$ws.Range("D45").Value = "'36.56"
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").Value = "'2.311.36"
$ws.Range("E47").Value = '  +1.79%  '
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = "'0.971"
$ws.Range("E49").Value = '  +2.46%  '
$ws.Range("D50").Value = "'20.70"
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("E51").Value = '  +1.17%  '
